$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: fix Objetivos (PT) content in B10/C10 ---
$ws.Range("B10").Value = "Geral: Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental. Executar procedimentos de análises volumétricas bem como interpretar, avaliar e criticar os resultados obtidos. Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Conhecer os métodos de preparação de amostras a serem analisadas. Compreender as bases teóricas da química analítica quantitativa de interesse ambiental. Executar procedimentos de análises químicas volumétricas (volumetria por neutralização, volumetria por precipitação, volumetria por oxi-redução e volumetria por complexação), análises ambientais de DBO, DQO, OD e turbidez, bem como interpretar, avaliar e criticar os resultados obtidos."
$ws.Range("C10").Value = "Geral: Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental. Executar procedimentos de análises volumétricas bem como interpretar, avaliar e criticar os resultados obtidos. Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Conhecer os métodos de preparação de amostras a serem analisadas. Compreender as bases teóricas da química analítica quantitativa de interesse ambiental. Executar procedimentos de análises químicas volumétricas (volumetria por neutralização, volumetria por precipitação, volumetria por oxi-redução e volumetria por complexação), análises ambientais de DBO, DQO, OD e turbidez, bem como interpretar, avaliar e criticar os resultados obtidos."

# --- Insert new row 13 for the teacher row (moves old 13.. down to 14..) ---
$ws.Rows.Item(13).Insert()

# Populate B13/C13 with the teacher string, using B10/C10 as format donors
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B13").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("C13").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("A13").Clear()

# --- Row 14 (was 13): Programa resumido -> new summary text ---
$ws.Range("B14").Value = "Fundamentos da análise titulométrica (preparação de amostras e cálculos). Análises titulométricas por precipitação, neutralização, complexação e oxirredução. Análises ambientais de DBO, DQO, OD e turbidez."
$ws.Range("C14").Value = "Fundamentos da análise titulométrica (preparação de amostras e cálculos). Análises titulométricas por precipitação, neutralização, complexação e oxirredução. Análises ambientais de DBO, DQO, OD e turbidez."

# --- Row 16 (was 15): Programa -> new full programa text ---
$ws.Range("B16").Value = "- Fundamentos de análise titulométrica e cálculos em análise titulométrica.- Titulometria de neutralização: fundamentos, indicadores de titulação, curvas de titulação ácido base.- Titulometria complexométrica: fundamentos, complexometria com EDTA.- Titulação de oxiredução: fundamentos e principais indicadores- Titulometria de precipitação: fundamentos, indicadores, argentimetria.- Análises quantitativas em solo, água, alimentos etc.- Análises ambientais: DBO, DQO, OD e turbidez"
$ws.Range("C16").Value = "- Fundamentos de análise titulométrica e cálculos em análise titulométrica.- Titulometria de neutralização: fundamentos, indicadores de titulação, curvas de titulação ácido base.- Titulometria complexométrica: fundamentos, complexometria com EDTA.- Titulação de oxiredução: fundamentos e principais indicadores- Titulometria de precipitação: fundamentos, indicadores, argentimetria.- Análises quantitativas em solo, água, alimentos etc.- Análises ambientais: DBO, DQO, OD e turbidez"

# --- Row 19 (was 18): Método -> method text ---
$ws.Range("B19").Value = "O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais."
$ws.Range("C19").Value = "O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais."

# --- Row 20 (was 19): Critério -> criterio text ---
$ws.Range("B20").Value = "Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."
$ws.Range("C20").Value = "Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."

# --- Row 21 (was 20): Norma de recuperação -> norma text ---
$ws.Range("B21").Value = "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 e estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."
$ws.Range("C21").Value = "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 e estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."

# --- Row 22 (was 21): Bibliografia -> new bibliography text ---
$ws.Range("B22").Value = "Harris, D.C. EXPLORANDO A QUÍMICA ANALÍTICA, 4ª edição, LTC, Rio de Janeiro – RJ, 2011Skoog, D.A., Holler, F.J. e Nieman, T.A., PRINCÍPIOS DE ANÁLISE INSTRUMENTAL, 5ª ed., Bookman, Porto Alegre, 2002.Mendham,J., Denney, R.C., Barnes, J.D. e Thomas, M., Vogel: ANÁLISE QUÍMICA QUANTITATIVA, 6ª ed., Livros Técnicos e Científicos, Rio de Janeiro -RJ, 2002."
$ws.Range("C22").Value = "Harris, D.C. EXPLORANDO A QUÍMICA ANALÍTICA, 4ª edição, LTC, Rio de Janeiro – RJ, 2011Skoog, D.A., Holler, F.J. e Nieman, T.A., PRINCÍPIOS DE ANÁLISE INSTRUMENTAL, 5ª ed., Bookman, Porto Alegre, 2002.Mendham,J., Denney, R.C., Barnes, J.D. e Thomas, M., Vogel: ANÁLISE QUÍMICA QUANTITATIVA, 6ª ed., Livros Técnicos e Científicos, Rio de Janeiro -RJ, 2002."

